# A new weekly price record was added to the "Arveja Verde" sheet.
# In Excel terms this is: select row 93, insert a new blank row above it
# (pushing the existing row 93 and everything below it down by one), then
# populate the newly-inserted row 93 with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93 - shifts rows 93:170 down to 94:171
$ws.Rows(93).Insert()

# Fill in the new row 93 with the new record's values
$ws.Cells.Item(93, 1).Value = 9
$ws.Cells.Item(93, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(93, 3).Value = "Metropolitana"
$ws.Cells.Item(93, 4).Value = 45118
$ws.Cells.Item(93, 5).Value = 13
$ws.Cells.Item(93, 6).Value = 100112022
$ws.Cells.Item(93, 7).Value = "Arveja Verde"
$ws.Cells.Item(93, 8).Value = "Perfection"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 43
$ws.Cells.Item(93, 11).Value = 36000
$ws.Cells.Item(93, 12).Value = 38000
$ws.Cells.Item(93, 13).Value = 37023
$ws.Cells.Item(93, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(93, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(93, 16).Value = 1481
$ws.Cells.Item(93, 17).Value = 25
$ws.Cells.Item(93, 18).Value = "Hortaliza"
